$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at E:F - shifts existing E.. content to G.. for every row
$ws.Columns("E:F").Insert()

# Restore/define column widths for the new columns (match column D's width of 18.5)
$ws.Columns("E:F").ColumnWidth = 17.666666666666668

# The insert left blank placeholder cells (carrying the row's style) in the
# other mini-tables (BD table row 5/6, Contract table row 8/9) that never had
# Static_Path/encap_vlan/Physical_Domain columns to begin with - clear them
# so they don't linger as empty-but-styled cells.
$ws.Range("E5:F5").Clear()
$ws.Range("E8:F8").Clear()

# ---- EPG table header row (row 1) ----
# E1 (Static_Path) was shifted to G1 by the column insert above; restore it.
# Note: write order matters for shared-string table layout -
# Physical_Domain, then encap_vlan (matches original authoring order)
$ws.Range("G1").Value = "Physical_Domain"
$ws.Range("F1").Value = "encap_vlan"
$ws.Range("E1").Value = "Static_Path"

# ---- EPG table data row (row 2) ----
$ws.Range("G2").Value = "E7_NETAPP"
$ws.Range("E2").Value = "e7_NETAPP-A_VPC"
$ws.Range("F2").Value = 800

# Selection / view state
$ws.Range("F2").Select()
